$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange {
    param($row1, $row2, $startCol, $endCol)
    $range1 = $ws.Range("$startCol$row1" + ":" + "$endCol$row1")
    $range2 = $ws.Range("$startCol$row2" + ":" + "$endCol$row2")
    $tmp = $range1.Value2
    $range1.Value2 = $range2.Value2
    $range2.Value2 = $tmp
}

Swap-RowRange 36 37 "B" "AD"
Swap-RowRange 99 100 "B" "AD"
